$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "JD_002"
$ws.Range("B3").Value = "Junior RPA Developer"
$ws.Range("C3").Value = "We are seeking a Junior RPA Developer to design, develop, and support automation solutions.`nCollaborate with teams to streamline business processes using RPA tools like UiPath or Automation Anywhere. Join Akkodis to grow your skills in a dynamic, tech-driven environment"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 4

$ws.Rows.Item(3).AutoFit()
